$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = -8.08
$ws.Range("C7").Value = -12.917
$ws.Range("B9").Value = 5.737
$ws.Range("C12").Value = -11.601
$ws.Range("C14").Value = -12.845
$ws.Range("D15").Value = -8.463000000000001
$ws.Range("B18").Value = 5.103
$ws.Range("B20").Value = 6.775999999999999
$ws.Range("C26").Value = -13.131
$ws.Range("B27").Value = 5.526000000000001
$ws.Range("C27").Value = -13.664
$ws.Range("C29").Value = -12.24
$ws.Range("D33").Value = -7.244000000000002
$ws.Range("B35").Value = 9.238999999999999
$ws.Range("D35").Value = -7.826000000000001
$ws.Range("C37").Value = -13.151
$ws.Range("C38").Value = -13.616
$ws.Range("D38").Value = -7.545
$ws.Range("D43").Value = -7.598999999999999
$ws.Range("D44").Value = -7.388
$ws.Range("D47").Value = -7.412999999999999
$ws.Range("C51").Value = -12.267
$ws.Range("D51").Value = -7.542999999999999
$ws.Range("C52").Value = -11.507
$ws.Range("C55").Value = -13.752
$ws.Range("D57").Value = -8.054
$ws.Range("D63").Value = -7.664999999999999
$ws.Range("B69").Value = 5.667
$ws.Range("C69").Value = -10.919
$ws.Range("C70").Value = -12.195
$ws.Range("D70").Value = -7.798
$ws.Range("B76").Value = 6.723999999999999
$ws.Range("B78").Value = 7.849999999999999
$ws.Range("C81").Value = -13.422
$ws.Range("B82").Value = 5.206999999999999
$ws.Range("B83").Value = 5.128
$ws.Range("C83").Value = -13.527
$ws.Range("D88").Value = -7.816999999999998
$ws.Range("B93").Value = 5.875999999999999
$ws.Range("D99").Value = -7.57
$ws.Range("C102").Value = -13.424
